# Generate Report for Handback
# Update the timestamps recorded on the "Overview", "zh-cn" and "de-de"
# sheets to reflect a fresh handback-status report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
# "Latest HO Xliff Generate Date" for 96f2d2a7-...md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 12:50:55"

# --- zh-cn sheet --------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the
# 96f2d2a7-...257b461171abb21adb83dbeab217cbee7f5e30dd.zh-cn.xlf row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-01 12:50:51"
$wsZhCn.Range("K3").Value = "2016-09-01 12:51:24"

# --- de-de sheet --------------------------------------------------------
# "Correspond Handoff Datetime" (shares text with Overview!G3) and
# "Correspond Handback DateTime" for the
# 96f2d2a7-...257b461171abb21adb83dbeab217cbee7f5e30dd.de-de.xlf row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-01 12:50:55"
$wsDeDe.Range("K3").Value = "2016-09-01 12:51:31"
